$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the headers in H1/I1: SamplePortion now comes before Result
$ws.Range("H1").Value = "SamplePortion"
$ws.Range("I1").Value = "Result"

# Update the #float type markers (row 2) with the unit annotation
$ws.Range("H2").Value = "#float,  unit:mg"
$ws.Range("I2").Value = "#float,  unit:mg"

# Add a new row 3 with French enum/description values for each column
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = "#PriseEssai"
$ws.Range("I3").Value = "#Resultat"
$ws.Range("J3").Value = "#NuméroLotReactif"
